$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header month / dates ---
# Dates are forced to text (NumberFormat "@") first so Excel's
# autodetection does not turn "2024.3.10" style strings into date serials.
$ws.Range("A1").Value = "3月"
$ws.Range("B1:F1").NumberFormat = "@"
$ws.Range("B1").Value = "2024.3.10"
$ws.Range("C1").Value = "2024.3.11"
$ws.Range("D1").Value = "2024.3.12"
$ws.Range("E1").Value = "2024.3.13"
$ws.Range("F1").Value = "2024.3.14"

# --- Row 2: 作业 / 完成作业 (C2 becomes blank) ---
$ws.Range("C2").Value = ""

# --- Row 3: 读书 / 35分钟 (C3 becomes blank) ---
$ws.Range("C3").Value = ""

# --- Row 4: 学习 row (reshuffled categories) ---
$ws.Range("B4").Value = "拼音"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "汉字"
$ws.Range("E4").Value = "拼音"
$ws.Range("F4").Value = "数学"

# --- Row 5: homework detail text (C5 becomes blank) ---
$homework = "1. 学习一页书`n2. 一套卷子`n3. 1页练字`n4. 拼音考核`n5. 英文`n6. 乘法表"
$ws.Range("B5").Value = $homework
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = $homework
$ws.Range("E5").Value = $homework
$ws.Range("F5").Value = $homework

# --- Row 6: 锻炼 / 跳绳 700 (C6 becomes blank) ---
$ws.Range("B6").Value = "跳绳 700"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "跳绳 700"
$ws.Range("E6").Value = "跳绳 700"
$ws.Range("F6").Value = "跳绳 700"

# --- Row 7: 积分 / star rating simplified to a single star (C7 becomes blank) ---
$ws.Range("B7").Value = "☆"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "☆"
$ws.Range("E7").Value = "☆"
$ws.Range("F7").Value = "☆"

# --- Selection moves from E5 to B4 ---
$ws.Range("B4").Select()

# --- Window height of the book view (best effort; engine may not persist this) ---
$win = $wb.Windows.Item(1)
$win.Height = 1192
